# Refresh crypto price/volume figures pulled by the GitHub Actions scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.402.27'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.111.03'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'" + '334.48'
$ws.Range("E5").Value = '  +2.83%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = "'" + '0.5232'
$ws.Range("E7").Value = '  +1.54%  '
$ws.Range("D8").Value = "'" + '0.4570'
$ws.Range("E8").Value = '  +6.39%  '
$ws.Range("D9").Value = "'" + '53.15'
$ws.Range("E9").Value = '  +17.02%  '
$ws.Range("D10").Value = "'" + '0.08926'
$ws.Range("E10").Value = '  +3.02%  '
$ws.Range("D11").Value = "'" + '1.179'
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("D12").Value = "'" + '24.47'
$ws.Range("E12").Value = '  +2.24%  '
$ws.Range("D13").Value = '2.098.48'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").Value = "'" + '6.819'
$ws.Range("E14").Value = '  +3.68%  '
$ws.Range("D15").Value = "'" + '8.013'
$ws.Range("E15").Value = '  +5.60%  '
$ws.Range("D16").Value = "'" + '96.60'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").Value = "'" + '1.003'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("E18").Value = '  +2.16%  '
$ws.Range("D19").Value = "'" + '0.06636'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").Value = "'" + '19.29'
$ws.Range("E20").Value = '  +3.73%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").Value = "'" + '6.369'
$ws.Range("D23").Value = '30.466.63'
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("D24").Value = "'" + '12.45'
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").Value = "'" + '2.364'
$ws.Range("E25").Value = '  +4.12%  '
$ws.Range("D26").Value = '2.345.28'
$ws.Range("E26").Value = '  +1.68%  '
$ws.Range("D27").Value = "'" + '22.40'
$ws.Range("E27").Value = '  +1.89%  '
$ws.Range("D28").Value = "'" + '2.577'
$ws.Range("E28").Value = '  +4.02%  '
$ws.Range("D29").Value = "'" + '163.69'
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").Value = "'" + '132.78'
$ws.Range("E30").Value = '  +2.31%  '
$ws.Range("E31").Value = '  +7.11%  '
$ws.Range("D32").Value = "'" + '1.720'
$ws.Range("E32").Value = '  +16.65%  '
$ws.Range("D33").Value = "'" + '0.1074'
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("D34").Value = "'" + '6.398'
$ws.Range("E34").Value = '  +6.26%  '
$ws.Range("D35").Value = "'" + '3.923'
$ws.Range("E35").Value = '  +2.35%  '
$ws.Range("D36").Value = "'" + '10.52'
$ws.Range("E36").Value = '  +11.03%  '
$ws.Range("D37").Value = "'" + '0.02590'
$ws.Range("E37").Value = '  +2.24%  '
$ws.Range("D38").Value = "'" + '0.06836'
$ws.Range("E38").Value = '  +4.57%  '
$ws.Range("D39").Value = "'" + '5.575'
$ws.Range("E39").Value = '  +3.93%  '
$ws.Range("D40").Value = "'" + '12.79'
$ws.Range("E40").Value = '  +3.57%  '
$ws.Range("D41").Value = "'" + '0.2301'
$ws.Range("E41").Value = '  +4.00%  '
$ws.Range("D42").Value = "'" + '0.6927'
$ws.Range("E42").Value = '  +5.16%  '
$ws.Range("D43").Value = "'" + '1.246'
$ws.Range("E43").Value = '  +1.61%  '
$ws.Range("D44").Value = "'" + '2.358'
$ws.Range("E44").Value = '  +8.57%  '
$ws.Range("D45").Value = "'" + '1.002'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").Value = "'" + '14.08'
$ws.Range("E46").Value = '  +1.58%  '
$ws.Range("D47").Value = "'" + '0.6389'
$ws.Range("E47").Value = '  +2.54%  '
$ws.Range("D48").Value = "'" + '3.660'
$ws.Range("E48").Value = '  +2.02%  '
$ws.Range("E49").Value = '  +25.98%  '
$ws.Range("D50").Value = "'" + '1.248'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("D51").Value = "'" + '0.3450'
$ws.Range("E51").Value = '  +29.16%  '
